$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "113766"
$ws.Range("C3").Value = "1065884"
$ws.Range("C4").Value = "1065909"
$ws.Range("C5").Value = "1071384"
$ws.Range("C6").Value = "1071902"
$ws.Range("C7").Value = "1071903"
$ws.Range("C8").Value = "1072080"
$ws.Range("C9").Value = "1072137"
$ws.Range("D2").Value = "174531"
$ws.Range("D3").Value = "174536"
$ws.Range("D4").Value = "367755"
$ws.Range("D5").Value = "1558064"
$ws.Range("D7").Value = "1068318"
$ws.Range("D8").Value = "1068327"
$ws.Range("D6").Value = "1068330"
$ws.Range("D9").Value = "135378"
$ws.Range("D10").Value = "135382"
$ws.Range("D11").Value = "135385"
$ws.Range("E2").Value = "1114394"
$ws.Range("E3").Value = "1532585"
$ws.Range("E4").Value = "1535683"
$ws.Range("E5").Value = "1535695"
$ws.Range("E6").Value = "1537215"
$ws.Range("E7").Value = "1566686"
$ws.Range("E8").Value = "1566802"
$ws.Range("E9").Value = "1566929"
$ws.Range("E10").Value = "1567159"
$ws.Range("E11").Value = "283088"
$ws.Range("E12").Value = "355307"
$ws.Range("E13").Value = "355319"
$ws.Range("E14").Value = "355574"
$ws.Range("E15").Value = "356165"
$ws.Range("E16").Value = "357035"
$ws.Range("E17").Value = "357038"
$ws.Range("E18").Value = "357041"
$ws.Range("E19").Value = "357044"
$ws.Range("E20").Value = "364831"
$ws.Range("E21").Value = "364834"
$ws.Range("E22").Value = "375935"
$ws.Range("E23").Value = "96889"
$ws.Range("E24").Value = "1665391"

$ws.Range("E25").Select() | Out-Null
